$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152, shifting existing rows 152:178 down to 153:179
$ws.Rows("152:152").Insert()

# Populate the newly inserted row 152 with the new data record
$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("D152").Value = 44522
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = 100112044
$ws.Range("G152").Value = "Perejil"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 60
$ws.Range("K152").Value = 5000
$ws.Range("L152").Value = 5000
$ws.Range("M152").Value = 5000
$ws.Range("N152").Value = "$/docena de atados (3 kilos)"
$ws.Range("O152").Value = "Región Metropolitana"
$ws.Range("P152").Value = 1667
$ws.Range("Q152").Value = 3
$ws.Range("R152").Value = "Hortaliza"
